$wb = $excel.ActiveWorkbook

# --- Sheet: Restricciones_del_follower (3rd sheet) ---
# Rows 2-5 (A:F) get new Expression / Function_Evaluation / Restriction_Set_Type /
# Lambda_value / Beta_value / Gamma_value content. All columns in this sheet are
# stored as text in the workbook (even the numeric-looking ones), so force the
# "@" (Text) number format on the numeric-looking columns (B, D, E, F) before
# writing so Excel doesn't silently coerce them to real numbers.
$wsFollower = $wb.Worksheets.Item(3)

$wsFollower.Range("B2:B5").NumberFormat = "@"
$wsFollower.Range("D2:F5").NumberFormat = "@"

$wsFollower.Range("A2").Value = "-534.0955555555555 - 2x_1 + 6.2333333333333325y_1 + 0.744444444444444y_2"
$wsFollower.Range("B2").Value = "536.5955555555555"
$wsFollower.Range("C2").Value = "J_0_L0_v"
$wsFollower.Range("D2").Value = "0.92"
$wsFollower.Range("E2").Value = "4.4"
$wsFollower.Range("F2").Value = "3.5"

$wsFollower.Range("A3").Value = "-853.0811111111111 + x_1 - 3x_2 + 8.366666666666667y_1 + 3.7888888888888888y_2"
$wsFollower.Range("B3").Value = "851.0811111111111"
$wsFollower.Range("C3").Value = "J_0_L0_v"
$wsFollower.Range("D3").Value = "0.36"
$wsFollower.Range("E3").Value = "10.0"
$wsFollower.Range("F3").Value = "5.2"

$wsFollower.Range("A4").Value = "10.159999999999998 - 0.09999999999999998y_1 + 0.3y_2"
$wsFollower.Range("B4").Value = "-10.159999999999998"
$wsFollower.Range("C4").Value = "J_0_LP_v"
$wsFollower.Range("D4").Value = "0.49"
$wsFollower.Range("E4").Value = "2.1"
$wsFollower.Range("F4").Value = "5.300000000000001"

$wsFollower.Range("A5").Value = "-883.4722222222221 + 8.633333333333331y_1 + 1.8777777777777769y_2"
$wsFollower.Range("B5").Value = "882.7022222222221"
$wsFollower.Range("C5").Value = "J_Ne_L0_v"
$wsFollower.Range("D5").Value = "0.75"
$wsFollower.Range("E5").Value = "8.0"
$wsFollower.Range("F5").Value = "5.0"

# --- Sheet: Punto_modificado (4th sheet) ---
# All values here are stored as text too.
$wsPunto = $wb.Worksheets.Item(4)
$wsPunto.Range("A2:D2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "51.550000000000004"
$wsPunto.Range("B2").Value = "18.099999999999998"
$wsPunto.Range("C2").Value = "102.2"
$wsPunto.Range("D2").Value = "0.2"

# --- Sheet: Vector_bf (5th sheet) ---
$wsBf = $wb.Worksheets.Item(5)
$wsBf.Range("A2:A3").NumberFormat = "@"
$wsBf.Range("A2").Value = "-11.172666666666665"
$wsBf.Range("A3").Value = "-4.604222222222221"

# --- Sheet: Vector_BF (6th sheet) ---
$wsBF = $wb.Worksheets.Item(6)
$wsBF.Range("A2:A5").NumberFormat = "@"
$wsBF.Range("A2").Value = "0.8000000000000007"
$wsBF.Range("A3").Value = "29.0"
$wsBF.Range("A4").Value = "-180.45"
$wsBF.Range("A5").Value = "-56.81666666666666"

# --- Sheet: Vector_Alpha (7th sheet) ---
# These cells are genuine numbers (not shared strings) in the workbook.
$wsAlpha = $wb.Worksheets.Item(7)
$wsAlpha.Range("A2").Value = 0.54
$wsAlpha.Range("A3").Value = 0.18
